# Apply changes described by the diff for data/PAL.MI.xlsx
#
# Row 109: update the timestamp (A109) and the "open" value (E109)
# Row 110: new row appended with a full set of OHLC data for 2024-07-22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 109 -------------------------------------------------
$ws.Range("A109").Value = 45492.2916666667
$ws.Range("E109").Value = 6.01999998092651

# --- Append new row 110 ------------------------------------------------------
# A110 needs the same date/time style as column A elsewhere (style index 1,
# numFmt "yyyy-mm-dd hh:mm:ss"). Copying the existing A109 cell brings that
# style across without introducing a brand-new style entry; the value is
# then overwritten with the correct timestamp.
$ws.Range("A109").Copy($ws.Range("A110"))
$ws.Range("A110").Value = 45495.4300462963

$ws.Range("B110").Value = 600
$ws.Range("C110").Value = 5.96000003814697
$ws.Range("D110").Value = 5.96000003814697
$ws.Range("E110").Value = 5.96000003814697
$ws.Range("F110").Value = 5.96000003814697

# G110 stores a numeric-looking value as text (shared string), matching how
# the "adj_close" column is written elsewhere in the sheet. Forcing the
# cell to Text format keeps it from being auto-converted to a number, and
# resetting the style back to Normal afterwards avoids leaving a stray
# number format applied to the cell.
$ws.Range("G110").NumberFormat = "@"
$ws.Range("G110").Value = "5.96000003814697"
$ws.Range("G110").Style = "Normal"

$ws.Range("H110").Value = "PAL.MI"
